# REVER_DailyTracker_MONISHA.xlsx — "Add files via upload"
# Target sheet is the one active when the file was saved (MAR-2021 , tab 7 / activeTab=6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) The whole date column (B2:B29) was off by exactly 10 years (3652 days) —
#    shift every date back to the correct MAR-2021 dates.
# ---------------------------------------------------------------------------
$dates = @{
    2  = 44256
    3  = 44257
    4  = 44258
    5  = 44259
    6  = 44260
    7  = 44261
    8  = 44262
    9  = 44263
    10 = 44264
    11 = 44265
    12 = 44266
    13 = 44267
    14 = 44268
    15 = 44269
    16 = 44270
    17 = 44271
    18 = 44272
    19 = 44273
    20 = 44274
    21 = 44275
    22 = 44276
    23 = 44277
    24 = 44278
    25 = 44279
    26 = 44280
    27 = 44281
    28 = 44282
    29 = 44283
}
foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 2).Value = $dates[$row]
}

# ---------------------------------------------------------------------------
# 2) Row 23 (day 22): task finished — % of completion goes to 100% and the
#    status flips from WIP to Completed (green "Completed" fill, style 19).
# ---------------------------------------------------------------------------
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "Completed"
$ws.Range("F4").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Rows 24 and 25 were blank placeholders — fill in the new tasks that were
#    logged for these two days.
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = "Hayaai -B2B"
$ws.Range("D24").Value = "Modification screen"
$ws.Range("E24").Value = 0.9
$ws.Range("F24").Value = "WIP"

$ws.Range("C25").Value = "Qmvar-2.0"
$ws.Range("D25").Value = "design issues checked"
$ws.Range("E25").Value = 0.9
$ws.Range("F25").Value = "WIP"

$ws.Range("F2").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("F25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Scroll/selection state of the sheet when it was saved.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F33").Select() | Out-Null
